$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Pattern" / "Pattern Type" header columns (M, N) next to the
# existing "Unique" column (L) for both property-type tables on the sheet
# (row 4 header block and row 11 header block), copying the style of the
# "Unique" header cell so the new cells look consistent.

$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)
$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

$ws.Range("L11").Copy()
$ws.Range("M11:N11").PasteSpecial(-4122)
$ws.Range("M11").Value = "Pattern"
$ws.Range("N11").Value = "Pattern Type"

$excel.CutCopyMode = 0

$ws.Range("M4:N4").Select()
